$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grant's normalized cell counts, replacing the raw Incucyte counts in B1:G23.
# Build a 23-row x 6-col array (columns B..G) and write it in a single block assign
# so the saved sheet matches an Excel paste over that range.
$normalized = New-Object 'object[,]' 23,6
$normalized[0,0] = 1588.7096378146102; $normalized[0,1] = 1613.2157188650124; $normalized[0,2] = 1627.5700325732901; $normalized[0,3] = 1551.1751212478698; $normalized[0,4] = 1579.9135128270254; $normalized[0,5] = 1354.7780514504655
$normalized[1,0] = 1598.0052793124619; $normalized[1,1] = 1719.7176777870118; $normalized[1,2] = 1738.8360477741587; $normalized[1,3] = 1629.7905361122032; $normalized[1,4] = 1568.3511557023112; $normalized[1,5] = 1423.2178434592229
$normalized[2,0] = 1641.1032535297729; $normalized[2,1] = 1737.0552059836164; $normalized[2,2] = 1802.1845819761131; $normalized[2,3] = 1654.3050203172104; $normalized[2,4] = 1652.5911861823722; $normalized[2,5] = 1447.9551176792559
$normalized[3,0] = 1723.0739103744629; $normalized[3,1] = 1806.4053187700345; $normalized[3,2] = 1850.914223669924; $normalized[3,3] = 1705.8699698518808; $normalized[3,4] = 1674.8900177800354; $normalized[3,5] = 1472.6923918992886
$normalized[4,0] = 1810.9599754450585; $normalized[4,1] = 1884.0113973643595; $normalized[4,2] = 1950.8099891422369; $normalized[4,3] = 1799.7012714641496; $normalized[4,4] = 1781.4288798577595; $normalized[4,5] = 1546.0796387520527
$normalized[5,0] = 1918.2823818293434; $normalized[5,1] = 1949.2335272468242; $normalized[5,2] = 2023.9044516829533; $normalized[5,3] = 1834.3596801677807; $normalized[5,4] = 1835.1112522225042; $normalized[5,5] = 1626.8880678708267
$normalized[6,0] = 1978.2815224063845; $normalized[6,1] = 2054.0842930072422; $normalized[6,2] = 2094.5624321389796; $normalized[6,3] = 1973.838642023856; $normalized[6,4] = 1967.2524765049527; $normalized[6,5] = 1704.398193760263
$normalized[7,0] = 2067.0126457949664; $normalized[7,1] = 2117.6552297281255; $normalized[7,2] = 2147.3528773072749; $normalized[7,3] = 2016.1049941014546; $normalized[7,4] = 1994.5066040132076; $normalized[7,5] = 1781.9083196496993
$normalized[8,0] = 2185.3208103130755; $normalized[8,1] = 2213.4244330998458; $normalized[8,2] = 2273.2377850162866; $normalized[8,3] = 2128.5334906278672; $normalized[8,4] = 2076.268986537973; $normalized[8,5] = 1871.7870826491519
$normalized[9,0] = 2257.9958256599143; $normalized[9,1] = 2306.716846729194; $normalized[9,2] = 2369.8849077090122; $normalized[9,3] = 2218.1381570323761; $normalized[9,4] = 2191.0666751333501; $normalized[9,5] = 1971.5607553366176
$normalized[10,0] = 2459.1197053407; $normalized[10,1] = 2541.1862756737505; $normalized[10,2] = 2524.1954397394138; $normalized[10,3] = 2399.0381439244984; $normalized[10,4] = 2387.626746253492; $normalized[10,5] = 2199.14367816092
$normalized[11,0] = 2609.5400859422962; $normalized[11,1] = 2717.038347382168; $normalized[11,2] = 2709.3680781758958; $normalized[11,3] = 2551.1970114038536; $normalized[11,4] = 2506.553848107696; $normalized[11,5] = 2291.4961685823755
$normalized[12,0] = 2766.720933087784; $normalized[12,1] = 2887.9368396058412; $normalized[12,2] = 2844.99891422367; $normalized[12,3] = 2704.2012059247604; $normalized[12,4] = 2664.2974345948687; $normalized[12,5] = 2420.9545703338808
$normalized[13,0] = 2975.4503376304483; $normalized[13,1] = 3081.952036091654; $normalized[13,2] = 2995.2486427795875; $normalized[13,3] = 2898.6264254817143; $normalized[13,4] = 2808.0010160020315; $normalized[13,5] = 2592.4663382594422
$normalized[14,0] = 3175.7291589932474; $normalized[14,1] = 3220.6522616644902; $normalized[14,2] = 3164.9902280130295; $normalized[14,3] = 3027.961462839166; $normalized[14,4] = 2971.5257810515618; $normalized[14,5] = 2745.0128626163114
$normalized[15,0] = 3230.6579496623699; $normalized[15,1] = 3332.9333966520244; $normalized[15,2] = 3270.5711183496201; $normalized[15,3] = 3074.4544501245246; $normalized[15,4] = 3156.5234950469894; $normalized[15,5] = 2854.6814449917902
$normalized[16,0] = 3316.0088397790059; $normalized[16,1] = 3398.1555265344891; $normalized[16,2] = 3301.4332247557004; $normalized[16,3] = 3154.7605190719619; $normalized[16,4] = 3170.5635001269998; $normalized[16,5] = 2883.5415982484951
$normalized[17,0] = 3339.6704726826279; $normalized[17,1] = 3389.0739641457913; $normalized[17,2] = 3309.5548317046691; $normalized[17,3] = 3154.7605190719619; $normalized[17,4] = 3191.2105664211322; $normalized[17,5] = 2918.1737821565412
$normalized[18,0] = 3368.4024554941684; $normalized[18,1] = 3455.9472871898374; $normalized[18,2] = 3373.7155266015202; $normalized[18,3] = 3198.7175252326642; $normalized[18,4] = 3201.1211582423161; $normalized[18,5] = 2938.7881773399017
$normalized[19,0] = 3400.5146715776555; $normalized[19,1] = 3469.9824290632791; $normalized[19,2] = 3374.5276872964173; $normalized[19,3] = 3215.624066063704; $normalized[19,4] = 3196.1658623317244; $normalized[19,5] = 2949.5076628352494
$normalized[20,0] = 3385.3036218538982; $normalized[20,1] = 3446.86572480114; $normalized[20,2] = 3363.1574375678611; $normalized[20,3] = 3224.9226635207756; $normalized[20,4] = 3202.7729235458464; $normalized[20,5] = 2964.3500273672689
$normalized[21,0] = 3386.9937384898712; $normalized[21,1] = 3450.9937077050931; $normalized[21,2] = 3340.4169381107495; $normalized[21,3] = 3218.1600471883598; $normalized[21,4] = 3187.0811531623058; $normalized[21,5] = 3000.6313628899838
$normalized[22,0] = 3423.3312461632909; $normalized[22,1] = 3468.3312359016977; $normalized[22,2] = 3357.4723127035832; $normalized[22,3] = 3222.3866823961198; $normalized[22,4] = 3230.0270510541018; $normalized[22,5] = 3008.8771209633283

$ws.Range("B1:G23").Value = $normalized

# Mirror the saved selection state (B1:G23 highlighted, active cell B1).
$null = $ws.Range("B1:G23").Select()
